$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45204
# (2023-10-05) to 45205 (2023-10-06) for every data row (rows 2-387).
$ws.Range("C2:C387").Value = 45205
